# Apply cryptos list refresh (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking string into a cell while keeping it
# stored as text (matches the original inline-string cells), by briefly
# switching the cell to Text format for the assignment and then restoring
# the General format so the cell keeps its original look/format.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
}

# --- Rows 16/17: Polkadot and WrappedBTC swapped ranking positions ---
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D16") "6.29"
$ws.Range("E16").Value = "  +7.95%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws.Range("D17") "57.681.53"
$ws.Range("E17").Value = "  +2.12%  "

# --- Price / Volume(1h) refresh for remaining rows ---
Set-TextValue $ws.Range("D2") "57.632.23"
$ws.Range("E2").Value = "  +2.10%  "
Set-TextValue $ws.Range("D3") "3.030.42"
$ws.Range("E3").Value = "  +1.41%  "
$ws.Range("E4").Value = "  +0.03%  "
Set-TextValue $ws.Range("D5") "512.62"
$ws.Range("E5").Value = "  +1.39%  "
Set-TextValue $ws.Range("D6") "139.95"
$ws.Range("E6").Value = "  +2.99%  "
$ws.Range("E7").Value = "  +0.02%  "
Set-TextValue $ws.Range("D8") "0.441"
$ws.Range("E8").Value = "  +2.30%  "
Set-TextValue $ws.Range("D9") "7.56"
$ws.Range("E9").Value = "  +0.72%  "
$ws.Range("E10").Value = "  +2.61%  "
Set-TextValue $ws.Range("D11") "0.368"
$ws.Range("E11").Value = "  +4.75%  "
Set-TextValue $ws.Range("D12") "3.553.49"
$ws.Range("E12").Value = "  +1.55%  "
$ws.Range("E13").Value = "  +1.77%  "
Set-TextValue $ws.Range("D14") "26.72"
$ws.Range("E14").Value = "  +4.30%  "
$ws.Range("E15").Value = "  +9.00%  "
Set-TextValue $ws.Range("D18") "3.035.98"
$ws.Range("E18").Value = "  +1.52%  "
Set-TextValue $ws.Range("D19") "12.86"
$ws.Range("E19").Value = "  +3.59%  "
Set-TextValue $ws.Range("D20") "8.04"
$ws.Range("E20").Value = "  +3.16%  "
Set-TextValue $ws.Range("D21") "333.16"
$ws.Range("E21").Value = "  +2.60%  "
$ws.Range("E22").Value = "  -0.09%  "
Set-TextValue $ws.Range("D23") "0.500"
$ws.Range("E23").Value = "  +5.53%  "
Set-TextValue $ws.Range("D24") "64.77"
$ws.Range("E24").Value = "  +3.95%  "
Set-TextValue $ws.Range("D25") "0.170"
$ws.Range("E25").Value = "  +3.61%  "
$ws.Range("E26").Value = "  +0.13%  "
Set-TextValue $ws.Range("D27") "0.0₃0935"
$ws.Range("E27").Value = "  +3.62%  "
Set-TextValue $ws.Range("D28") "6.82"
$ws.Range("E28").Value = "  +5.80%  "
Set-TextValue $ws.Range("D29") "7.51"
$ws.Range("E29").Value = "  +9.06%  "
$ws.Range("E30").Value = "  +3.39%  "
Set-TextValue $ws.Range("D31") "1.21"
$ws.Range("E31").Value = "  +0.61%  "
Set-TextValue $ws.Range("D32") "20.77"
$ws.Range("E32").Value = "  +1.39%  "
Set-TextValue $ws.Range("D33") "4.74"
$ws.Range("E33").Value = "  +6.30%  "
Set-TextValue $ws.Range("D34") "155.72"
$ws.Range("E34").Value = "  -0.80%  "
Set-TextValue $ws.Range("D35") "5.89"
$ws.Range("E35").Value = "  +5.78%  "
Set-TextValue $ws.Range("D36") "1.28"
$ws.Range("E36").Value = "  +2.31%  "
Set-TextValue $ws.Range("D37") "24.79"
$ws.Range("E37").Value = "  +6.49%  "
$ws.Range("E38").Value = "  +1.99%  "
Set-TextValue $ws.Range("D39") "3.068.00"
$ws.Range("E39").Value = "  +1.55%  "
Set-TextValue $ws.Range("D40") "37.51"
$ws.Range("E40").Value = "  +3.27%  "
$ws.Range("E41").Value = "  +8.59%  "
$ws.Range("E42").Value = "  +0.10%  "
Set-TextValue $ws.Range("D43") "2.319.10"
$ws.Range("E43").Value = "  +2.90%  "
Set-TextValue $ws.Range("D44") "0.656"
$ws.Range("E44").Value = "  +2.06%  "
$ws.Range("E45").Value = "  +1.93%  "
Set-TextValue $ws.Range("D46") "0.997"
$ws.Range("E46").Value = "  +1.30%  "
Set-TextValue $ws.Range("D47") "6.05"
$ws.Range("E47").Value = "  +5.30%  "
Set-TextValue $ws.Range("D48") "0.0241"
$ws.Range("E48").Value = "  +2.30%  "
Set-TextValue $ws.Range("D49") "19.63"
$ws.Range("E49").Value = "  +3.42%  "
Set-TextValue $ws.Range("D50") "1.86"
$ws.Range("E50").Value = "  -5.12%  "
Set-TextValue $ws.Range("D51") "0.0895"
$ws.Range("E51").Value = "  +2.97%  "
